$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = '28.703.93'
$ws.Range("E2").Value = '  +1.69%  '

# Row 3: D3, E3
$ws.Range("D3").Value = '1.572.87'
$ws.Range("E3").Value = '  -1.10%  '

# Row 5: D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.01'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6: E6
$ws.Range("E6").Value = '  +0.08%  '

# Row 7: E7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8: D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.69'
$ws.Range("E8").Value = '  +4.00%  '

# Row 9: D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.20'
$ws.Range("E9").Value = '  +0.50%  '

# Row 10: E10
$ws.Range("E10").Value = '  -1.49%  '

# Row 11: E11
$ws.Range("E11").Value = '  -1.26%  '

# Row 12: E12
$ws.Range("E12").Value = '  -0.11%  '

# Row 13: D13, E13
$ws.Range("D13").Value = '1.797.11'
$ws.Range("E13").Value = '  -1.14%  '

# Row 14: D14, E14
$ws.Range("D14").Value = '1.577.13'
$ws.Range("E14").Value = '  -0.81%  '

# Row 15: B15, C15, D15, E15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '28.706.85'
$ws.Range("E15").Value = '  +1.55%  '

# Row 16: B16, C16, D16, E16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.522'
$ws.Range("E16").Value = '  -1.67%  '

# Row 17: D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.68'
$ws.Range("E17").Value = '  -1.77%  '

# Row 18: D18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.43'

# Row 19: D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.69'
$ws.Range("E19").Value = '  +1.44%  '

# Row 20: E20
$ws.Range("E20").Value = '  -1.12%  '

# Row 21: E21
$ws.Range("E21").Value = '  -2.20%  '

# Row 22: D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.11%  '

# Row 23: D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.90'
$ws.Range("E23").Value = '  -5.10%  '

# Row 24: D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.17'
$ws.Range("E24").Value = '  -1.63%  '

# Row 25: D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.13'
$ws.Range("E25").Value = '  +9.66%  '

# Row 26: D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.87'
$ws.Range("E26").Value = '  +0.01%  '

# Row 27: E27
$ws.Range("E27").Value = '  -1.02%  '

# Row 28: D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.45'
$ws.Range("E28").Value = '  -1.97%  '

# Row 29: E29
$ws.Range("E29").Value = '  -2.71%  '

# Row 30: E30
$ws.Range("E30").Value = '  -0.06%  '

# Row 31: D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("E31").Value = '  +2.76%  '

# Row 32: E32
$ws.Range("E32").Value = '  -2.42%  '

# Row 33: E33
$ws.Range("E33").Value = '  -0.63%  '

# Row 34: E34
$ws.Range("E34").Value = '  -1.17%  '

# Row 35: D35, E35
$ws.Range("D35").Value = '1.389.03'
$ws.Range("E35").Value = '  -0.74%  '

# Row 36: E36
$ws.Range("E36").Value = '  +2.12%  '

# Row 37: E37
$ws.Range("E37").Value = '  -3.09%  '

# Row 38: E38
$ws.Range("E38").Value = '  +0.73%  '

# Row 39: D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.61'
$ws.Range("E39").Value = '  +2.20%  '

# Row 40: E40
$ws.Range("E40").Value = '  -0.24%  '

# Row 41: D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.526'
$ws.Range("E41").Value = '  -2.93%  '

# Row 42: B42, C42, D42, E42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.91'
$ws.Range("E42").Value = '  +1.76%  '

# Row 43: B43, C43, D43, E43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.05%  '

# Row 44: D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.795'
$ws.Range("E44").Value = '  -2.43%  '

# Row 45: E45
$ws.Range("E45").Value = '  +2.12%  '

# Row 46: D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.53'
$ws.Range("E46").Value = '  -1.37%  '

# Row 47: E47
$ws.Range("E47").Value = '  -1.56%  '

# Row 48: D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.25'
$ws.Range("E48").Value = '  -1.64%  '

# Row 49: D49, E49
$ws.Range("D49").Value = '1.710.01'
$ws.Range("E49").Value = '  -1.12%  '

# Row 50: D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.70'
$ws.Range("E50").Value = '  -0.79%  '

# Row 51: B51, C51, D51, E51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -0.57%  '
